$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (after "MP73010" in the title line).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Merge the ">>>  your stuff after this line >>>" paragraph's three runs
#    (and the gramStart/gramEnd proof-error marks) into a single run by
#    doing a Find/Replace across the full paragraph text.
$d.Content.Find.Execute(">>>  your stuff after this line >>>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, ">>>  your stuff after this line >>>", 2) | Out-Null

# 3. Replace "Ben changing things up!" with the Version Control paragraph text.
$dash = [char]0x2013
$versionText = "Version Control $dash this keeps track of the many versions of system components (e.g. source files, media objects, help files, etc.) that may exist during development."
$d.Content.Find.Execute("Ben changing things up!", $false, $false, $false, $false, $false, `
                         $true, 1, $false, $versionText, 2) | Out-Null

# 4. Turn the following empty paragraph into "Edited by Rocky" and put the
#    _GoBack bookmark at the end of that run (collapsed, after "Rocky").
#    A trailing sentinel character is used while placing the bookmark because
#    collapsed ranges positioned exactly at a paragraph's text end (just
#    before the paragraph mark) are mis-resolved; having a character after
#    the insertion point keeps it away from that boundary. The sentinel is
#    removed once the bookmark is anchored.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq [string][char]13) {
        $targetPara = $p
        break
    }
}
$targetPara.Range.InsertBefore("Edited by RockyX")

$found = $d.Content
$found.Find.Execute("Edited by RockyX") | Out-Null
$sentinelStart = $found.End - 1

$bm = $d.Range($sentinelStart, $sentinelStart)
$bm.Bookmarks.Add("_GoBack") | Out-Null

$sentinel = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinel.Text = ""
